$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> updated nombre_aides (C) and montant_total (D) values
# reflecting the 2020-09-22 data refresh for Fonds de solidarite volet 2 (regional x NAF)
$updates = @(
    ,@(3, "95", "333198.00")
    ,@(4, "9", "59000.00")
    ,@(6, "179", "461886.40")
    ,@(7, "529", "1338809.82")
    ,@(8, "464", "1058821.00")
    ,@(9, "926", "3908543.29")
    ,@(14, "126", "358245.00")
    ,@(15, "102", "245788.98")
    ,@(16, "26", "55593.58")
    ,@(17, "155", "845186.23")
    ,@(33, "19", "51500.00")
    ,@(34, "35", "113000.00")
    ,@(36, "175", "492408.00")
    ,@(38, "395", "1595050.18")
    ,@(44, "34", "88971.00")
    ,@(46, "60", "244138.95")
    ,@(47, "82", "191125.00")
    ,@(53, "286", "1295408.10")
    ,@(66, "77", "195000.00")
    ,@(67, "139", "489413.09")
    ,@(71, "45", "152724.00")
    ,@(74, "15", "55071.60")
    ,@(79, "101", "275987.00")
    ,@(80, "231", "637575.19")
    ,@(81, "35", "119024.78")
    ,@(82, "572", "2834230.05")
    ,@(85, "22", "70569.00")
    ,@(86, "85", "309474.36")
    ,@(122, "40", "170158.69")
    ,@(123, "68", "182877.00")
    ,@(124, "275", "795824.00")
    ,@(125, "225", "611728.39")
    ,@(126, "635", "3403293.37")
    ,@(127, "13", "45000.00")
    ,@(131, "52", "220279.76")
    ,@(134, "104", "539262.28")
    ,@(207, "707", "3030474.95")
    ,@(215, "140", "732021.14")
    ,@(217, "58", "151000.00")
    ,@(218, "69", "197529.05")
    ,@(219, "5", "13000.00")
    ,@(220, "155", "394020.07")
    ,@(221, "464", "1260364.16")
    ,@(222, "51", "144412.00")
    ,@(223, "880", "3498781.86")
    ,@(224, "26", "71518.00")
    ,@(226, "58", "136106.00")
    ,@(227, "184", "579305.17")
    ,@(228, "71", "233420.52")
    ,@(229, "90", "249892.59")
    ,@(230, "28", "63500.00")
    ,@(231, "155", "641455.74")
    ,@(232, "186", "432145.00")
    ,@(249, "97", "292338.00")
    ,@(250, "187", "480590.00")
    ,@(251, "624", "1654568.61")
    ,@(252, "211", "574949.11")
    ,@(253, "1137", "4433238.95")
    ,@(254, "54", "154503.19")
    ,@(256, "91", "218500.00")
    ,@(257, "216", "706297.64")
    ,@(258, "154", "528293.00")
    ,@(259, "111", "299582.92")
    ,@(260, "35", "80500.00")
    ,@(261, "171", "633874.29")
    ,@(262, "249", "580337.52")
)

foreach ($u in $updates) {
    $row = $u[0]
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $u[1]
    $cCell.ClearFormats()

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u[2]
    $dCell.ClearFormats()
}
